$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 24; $i++) {
    $ws.Cells.Item($i, 1).Value = 0.0
}
